$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Table header row, column 2: "StringToWord" -> "StringToWord Attribute Filter"
$cell2 = $t.Cell(2, 2)
$r2 = $cell2.Range
$r2.End = $r2.End - 1              # exclude the cell end-of-cell mark
$r2.InsertAfter(" Attribute Filter")

# --- Table header row, column 3: "Personal Tokenization" -> "Tokenization"
#     plus the "_GoBack" bookmark now sits at the start of this cell's paragraph
$cell3 = $t.Cell(2, 3)
$r3 = $cell3.Range
$r3.End = $r3.End - 1              # exclude the cell end-of-cell mark
$r3.Text = "Tokenization"

$cell3b = $t.Cell(2, 3)
$bmStart = $cell3b.Range.Start
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Body paragraph describing Test 1
$d.Content.Find.Execute(
    "Test 1 uses the provided",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Test 1 is computed on the provided",
    2)

$d.Content.Find.Execute(
    "file, and uses StringToWord attribute filter",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "file, and uses the StringToWord attribute filter",
    2)
